# Auto-generated: update Vacancy / pricing cells per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Studio: 10`n1 bed: 4`nPercent: 6.93%"
$ws.Range("D2").Value = "6/16: `$1920-2120`n(`$3.91-4.51)"
$ws.Range("E2").Value = "6/16: `$2375-2420`n(`$3.86-4.51)"
$ws.Range("B3").Value = "1 bed: 1`n2 bed: 2`nPercent: 1.88%"
$ws.Range("E3").Value = "6/16: `$2852`n(`$3.19)"
$ws.Range("F3").Value = "6/16: `$3600-3650`n(`$2.73-2.91)"
$ws.Range("B4").Value = "1 bed: 7`n2 bed: 3`nPercent: 4.42%"
$ws.Range("E4").Value = "6/16: `$2250-2600`n(`$3.54-4.01)"
$ws.Range("F4").Value = "6/16: `$2980-3500`n(`$3.39-3.74)"
$ws.Range("B5").Value = "1 bed: 11`n2 bed: 1`nPercent: 5.74%"
$ws.Range("C5").Value = "UP TO 8 WEEKS FREE Special on select units! Some restrictions May Apply. Please contact the Leasing Office for details."
$ws.Range("E5").Value = "6/16: `$2000-2400`n(`$3.46-4.14)"
$ws.Range("F5").Value = "6/16: `$3300`n(`$3.44)"
$ws.Range("B6").Value = "1 bed: 6`n2 bed: 8`nPercent: 11.02%"
$ws.Range("C6").Value = "Move-in Special Now Offering 2 Months FREE on a 15 Month Lease OR 1 Month Free on 13 Month Lease with approval credit."
$ws.Range("E6").Value = "6/16: `$2250-2300`n(`$3.47-3.87)"
$ws.Range("F6").Value = "6/16: `$2750-3050`n(`$3.01-3.34)"
$ws.Range("B7").Value = "Studio: 1`n1 bed: 1`nPercent: 11.76%"
$ws.Range("D7").Value = "6/16: `$2250`n(`$3.16)"
$ws.Range("E7").Value = "6/16: `$3900`n(`$3.21)"
$ws.Range("B8").Value = "Studio: 1`nPercent: 0.91%"
$ws.Range("D8").Value = "6/16: `$1950`n(`$3.62)"
$ws.Range("E8").Value = "n/a"
$ws.Range("B9").Value = "Studio: 5`n1 bed: 2`n3 bed: 1`nPercent: 2.38%"
$ws.Range("D9").Value = "6/16: `$2050-2200`n(`$2.99-3.55)"
$ws.Range("E9").Value = "6/16: `$2675-3460`n(`$3.25-3.64)"
$ws.Range("G9").Value = "6/16: `$5260`n(`$3.66)"
$ws.Range("B10").Value = "1 bed: 6`n2 bed: 1`nPercent: 3.07%"
$ws.Range("E10").Value = "6/16: `$2675-4475`n(`$3.42-4.25)"
$ws.Range("F10").Value = "6/16: `$3275`n(`$3.27)"
$ws.Range("B11").Value = "Studio: 10`n1 bed: 13`n2 bed: 4`nPercent: 12.05%"
$ws.Range("D11").Value = "6/16: `$2017-2289`n(`$3.65-4.51)"
$ws.Range("E11").Value = "6/16: `$2120-2895`n(`$2.84-3.69)"
$ws.Range("F11").Value = "6/16: `$3092-3822`n(`$3.04-3.7)"
$ws.Range("B12").Value = "1 bed: 7`nPercent: 3.13%"
$ws.Range("E12").Value = "6/16: `$2290-3515`n(`$2.69-3.54)"
$ws.Range("F12").Value = "n/a"
$ws.Range("B13").Value = "1 bed: 5`n2 bed: 3`nPercent: 6.11%"
$ws.Range("E13").Value = "6/16: `$2701-2933`n(`$3.35-3.71)"
$ws.Range("F13").Value = "6/16: `$3224-3786`n(`$2.69-3.15)"
$ws.Range("B14").Value = "Studio: 18`n1 bed: 4`n2 bed: 2`nPercent: 7.29%"
$ws.Range("D14").Value = "6/16: `$2099-2445`n(`$3.11-3.88)"
$ws.Range("E14").Value = "6/16: `$2400-2956`n(`$2.58-3.63)"
$ws.Range("F14").Value = "6/16: `$3300-3328`n(`$2.88-2.9)"
$ws.Range("B15").Value = "Studio: 2`n1 bed: 3`n2 bed: 2`nPercent: 3.07%"
$ws.Range("D15").Value = "6/16: `$2100-2156`n(`$3.56-3.93)"
$ws.Range("E15").Value = "6/16: `$2463-2760`n(`$3.55-3.77)"
$ws.Range("F15").Value = "6/16: `$3295-3466`n(`$3.3-3.54)"
$ws.Range("B16").Value = "1 bed: 1`n2 bed: 2`nPercent: 3.06%"
$ws.Range("E16").Value = "6/16: `$2631`n(`$3.01)"
$ws.Range("F16").Value = "6/16: `$3245-3676`n(`$2.92-2.96)"
$ws.Range("B17").Value = "Studio: 7`n1 bed: 1`nPercent: 7.34%"
$ws.Range("D17").Value = "6/16: `$1780-2075`n(`$3.49-4.13)"
$ws.Range("E17").Value = "6/16: `$2125`n(`$3.84)"
$ws.Range("B18").Value = "Studio: 7`n1 bed: 10`n2 bed: 12`n3 bed: 1`nPercent: 4.66%"
$ws.Range("D18").Value = "6/16: `$2505-3248`n(`$4.09-5.19)"
$ws.Range("E18").Value = "6/16: `$2765-3883`n(`$3.3-5.18)"
$ws.Range("F18").Value = "6/16: `$4097-18995`n(`$3.55-6.78)"
$ws.Range("G18").Value = "6/16: `$20995`n(`$5.59)"
$ws.Range("B19").Value = "Studio: 13`n1 bed: 14`n2 bed: 2`nPercent: 11.69%"
$ws.Range("D19").Value = "6/16: `$1955-2212`n(`$3.29-4.4)"
$ws.Range("E19").Value = "6/16: `$2380-2769`n(`$3.38-4.09)"
$ws.Range("F19").Value = "6/16: `$3310-3435`n(`$3.78-3.92)"
